$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.532.77'
$ws.Range("E2").Value = '  -2.31%  '
$ws.Range("D3").Value = '2.891.34'
$ws.Range("E3").Value = '  -2.09%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '566.40'
$ws.Range("E5").Value = '  -4.52%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.93'
$ws.Range("E6").Value = '  -3.68%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.501'
$ws.Range("E8").Value = '  -1.08%  '
$ws.Range("D9").Value = '2.890.09'
$ws.Range("E9").Value = '  -2.08%  '
$ws.Range("E10").Value = '  -2.51%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.146'
$ws.Range("E11").Value = '  -2.65%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.429'
$ws.Range("E12").Value = '  -2.56%  '
$ws.Range("E13").Value = '  -1.53%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '31.76'
$ws.Range("E14").Value = '  -3.03%  '
$ws.Range("E15").Value = '  -0.49%  '
$ws.Range("D16").Value = '3.368.16'
$ws.Range("E16").Value = '  -2.17%  '
$ws.Range("D17").Value = '61.498.82'
$ws.Range("E17").Value = '  -2.37%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.53'
$ws.Range("E18").Value = '  -2.13%  '
$ws.Range("D19").Value = '2.887.70'
$ws.Range("E19").Value = '  -2.21%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '432.46'
$ws.Range("E20").Value = '  -1.87%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.01'
$ws.Range("E21").Value = '  -3.54%  '
$ws.Range("E22").Value = '  -2.35%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.82'
$ws.Range("E23").Value = '  -2.75%  '
$ws.Range("E24").Value = '  -2.11%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.88'
$ws.Range("E25").Value = '  +0.57%  '
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.93'
$ws.Range("E27").Value = '  -11.78%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.01'
$ws.Range("E28").Value = '  -6.13%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0000106'
$ws.Range("E29").Value = '  +2.93%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.99'
$ws.Range("E30").Value = '  -4.63%  '
$ws.Range("E31").Value = '  -4.59%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.04'
$ws.Range("E32").Value = '  -8.93%  '
$ws.Range("E33").Value = '  +0.03%  '
$ws.Range("E34").Value = '  -2.21%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '25.51'
$ws.Range("E35").Value = '  -3.40%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.955'
$ws.Range("E36").Value = '  -3.68%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.36'
$ws.Range("E37").Value = '  -4.69%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '48.84'
$ws.Range("E38").Value = '  -1.85%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.93'
$ws.Range("E39").Value = '  -5.42%  '
$ws.Range("E40").Value = '  -9.11%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.19'
$ws.Range("E41").Value = '  -3.66%  '
$ws.Range("E42").Value = '  -3.73%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '39.60'
$ws.Range("E43").Value = '  +0.05%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.267'
$ws.Range("E44").Value = '  -4.72%  '
$ws.Range("D45").Value = '2.686.36'
$ws.Range("E45").Value = '  -0.85%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '132.62'
$ws.Range("E46").Value = '  -2.34%  '
$ws.Range("E47").Value = '  -1.43%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '341.73'
$ws.Range("E49").Value = '  -5.02%  '
$ws.Range("E50").Value = '  -2.02%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '21.49'
$ws.Range("E51").Value = '  -5.56%  '
